$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.615.43"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "1.843.15"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'315.29"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "'0.4228"
$ws.Range("E7").Value = "  -2.70%  "
$ws.Range("D8").Value = "'0.3634"
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("D9").Value = "'45.34"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "'0.07251"
$ws.Range("D11").Value = "'0.8881"
$ws.Range("D12").Value = "'20.60"
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("D13").Value = "1.849.73"
$ws.Range("E13").Value = "  -4.07%  "
$ws.Range("D14").Value = "'6.556"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").Value = "'5.326"
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "'0.06857"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "'78.88"
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("D19").Value = "'0.000008853"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "'15.43"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("D22").Value = "27.604.16"
$ws.Range("E22").Value = "  -1.09%  "
$ws.Range("D23").Value = "'4.969"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("E24").Value = "  -4.57%  "
$ws.Range("D25").Value = "2.075.34"
$ws.Range("E25").Value = "  +0.70%  "
$ws.Range("D26").Value = "'1.958"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("D27").Value = "'155.29"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("D28").Value = "'18.71"
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").Value = "'122.47"
$ws.Range("E29").Value = "  +8.20%  "
$ws.Range("D30").Value = "'5.236"
$ws.Range("E30").Value = "  -2.45%  "
$ws.Range("D31").Value = "'1.848"
$ws.Range("E31").Value = "  +7.16%  "
$ws.Range("D32").Value = "'0.08900"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("D33").Value = "'0.7765"
$ws.Range("E33").Value = "  -2.92%  "
$ws.Range("D34").Value = "'4.568"
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("D35").Value = "'2.921"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").Value = "'1.090"
$ws.Range("E36").Value = "  -6.50%  "
$ws.Range("D37").Value = "'0.9996"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "'0.01925"
$ws.Range("E40").Value = "  -2.09%  "
$ws.Range("D41").Value = "'2.788"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("D42").Value = "'6.857"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("D43").Value = "'0.5067"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("D45").Value = "'8.261"
$ws.Range("E45").Value = "  -5.28%  "
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'10.38"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "'0.4704"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("D49").Value = "'104.55"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "'0.9998"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'1.630"
$ws.Range("E51").Value = "  -2.46%  "
